$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=3.669332; H=11.007996; I=0.1142408795870128; J=0.1260733720492186; K=3; M=8.889220333333332; N=26.667661; O=0.5709282459015017; P=0.5709282459015017; Q=32.61750062415067; R=293.557505617356; S=0.0652233449928579; T=0.07197884915894778 }
    3  = @{ E=3; G=3.669332; H=11.007996; I=0.1142408795870128; J=0.1260733720492186; K=3; M=6.680547666666667; N=20.041643; O=0.4290717540984982; P=0.4290717540984982; Q=24.51314733082533; R=220.618325977428; S=0.04901753459415491; T=0.0540945228902708 }
    4  = @{ E=3; G=1.273458666666667; H=3.820376; I=0.03964782641573577; J=0.04375434773194917; K=3; M=8.889220333333332; N=26.667661; O=0.5709282459015017; P=0.5709282459015017; Q=11.32005467339289; R=101.880492060536; S=0.02263606398934325; T=0.02498059300116609 }
    5  = @{ E=3; G=1.273458666666667; H=3.820376; I=0.03964782641573577; J=0.04375434773194917; K=3; M=6.680547666666667; N=20.041643; O=0.4290717540984982; P=0.4290717540984982; Q=8.507401324196445; R=76.56661191776801; S=0.01701176242639252; T=0.01877375473078308 }
    6  = @{ E=3; G=8.210356666666668; H=24.63107; I=0.2556210142126945; J=0.2820969459000845; K=3; M=8.889220333333332; N=26.667661; O=0.5709282459015017; P=0.5709282459015017; Q=72.98366942525222; R=656.85302482727; S=0.1459412572600165; T=0.1610571144969061 }
    7  = @{ E=3; G=8.210356666666668; H=24.63107; I=0.2556210142126945; J=0.2820969459000845; K=3; M=6.680547666666667; N=20.041643; O=0.4290717540984982; P=0.4290717540984982; Q=54.84967907200112; R=493.64711164801; S=0.109679756952678; T=0.1210398314031784 }
    8  = @{ E=3; G=9.922544666666667; H=29.767634; I=0.3089282273889152; J=0.3409254505821921; K=3; M=8.889220333333332; N=26.667661; O=0.5709282459015017; P=0.5709282459015017; Q=88.20368580934155; R=793.833172284074; S=0.1763758509726136; T=0.1946439694840701 }
    9  = @{ E=3; G=9.922544666666667; H=29.767634; I=0.3089282273889152; J=0.3409254505821921; K=3; M=6.680547666666667; N=20.041643; O=0.4290717540984982; P=0.4290717540984982; Q=66.28803262029578; R=596.592293582662; S=0.1325523764163015; T=0.1462814810981221 }
    10 = @{ E=2; G=9.043563500000001; H=18.087127; I=0.2815620523956416; J=0.2071498837365554; K=3; M=8.889220333333332; N=26.667661; O=0.5709282459015017; P=0.5709282459015017; Q=80.39022854999116; R=482.3413712999471; S=0.1607517286866704; T=0.1182677197604116 }
    11 = @{ E=2; G=9.043563500000001; H=18.087127; I=0.2815620523956416; J=0.2071498837365554; K=3; M=6.680547666666667; N=20.041643; O=0.4290717540984982; P=0.4290717540984982; Q=60.41595703827684; R=362.4957422296611; S=0.1208103237089712; T=0.0888821639761438 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
